# Update column G ("K" - strikeouts) values for rows 2-37 with freshly
# regenerated data (replacing the old "Strike#" derived values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 1
    25 = 2
    26 = 2
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 1
    36 = 3
    37 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
